$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices, volume %, and two row swaps)
# Each cell is forced to Text format while writing so that values like
# '66.105.71', '0.622', '9.92' etc. are stored as text (matching the
# source data) rather than being auto-converted to numbers by Excel,
# then the explicit style is reset back to Normal so no stray number
# format / style index is left behind on the cell.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '66.105.71'
$ws.Range('D2').Style = "Normal"

$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +3.87%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.798.15'
$ws.Range('D3').Style = "Normal"

$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +6.49%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '428.79'
$ws.Range('D5').Style = "Normal"

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +7.69%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '139.61'
$ws.Range('D6').Style = "Normal"

$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +13.66%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.622'
$ws.Range('D7').Style = "Normal"

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +5.95%  '
$ws.Range('E7').Style = "Normal"

$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +8.67%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0000316'
$ws.Range('D11').Style = "Normal"

$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -3.83%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '42.87'
$ws.Range('D12').Style = "Normal"

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +9.97%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '10.56'
$ws.Range('D13').Style = "Normal"

$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +14.60%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.394.55'
$ws.Range('D14').Style = "Normal"

$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +6.57%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '15.21'
$ws.Range('D15').Style = "Normal"

$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +9.42%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.23%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '3.785.33'
$ws.Range('D17').Style = "Normal"

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +6.07%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +6.90%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +10.34%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '66.246.60'
$ws.Range('D20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '407.62'
$ws.Range('D21').Style = "Normal"

$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +3.16%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '15.13'
$ws.Range('D22').Style = "Normal"

$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +8.79%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +10.40%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '85.26'
$ws.Range('D24').Style = "Normal"

$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +3.89%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '36.75'
$ws.Range('D25').Style = "Normal"

$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +7.98%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('B26').NumberFormat = "@"
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('B26').Style = "Normal"

$ws.Range('C26').NumberFormat = "@"
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('C26').Style = "Normal"

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '3.31'
$ws.Range('D26').Style = "Normal"

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +9.67%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('B27').NumberFormat = "@"
$ws.Range('B27').Value = 'RenderToken'
$ws.Range('B27').Style = "Normal"

$ws.Range('C27').NumberFormat = "@"
$ws.Range('C27').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('C27').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.79'
$ws.Range('D27').Style = "Normal"

$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +42.81%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.92'
$ws.Range('D28').Style = "Normal"

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +13.26%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.41'
$ws.Range('D29').Style = "Normal"

$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.27%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '13.89'
$ws.Range('D30').Style = "Normal"

$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +16.22%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '712.43'
$ws.Range('D31').Style = "Normal"

$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +6.23%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.133'
$ws.Range('D32').Style = "Normal"

$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +18.29%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  +7.31%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '41.17'
$ws.Range('D34').Style = "Normal"

$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +12.45%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  +40.70%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +1.27%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '56.53'
$ws.Range('D38').Style = "Normal"

$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +4.61%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.0476'
$ws.Range('D39').Style = "Normal"

$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +9.10%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.65'
$ws.Range('D40').Style = "Normal"

$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +45.61%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.0₃0687'
$ws.Range('D41').Style = "Normal"

$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +5.39%  '
$ws.Range('E41').Style = "Normal"

$ws.Range('B42').NumberFormat = "@"
$ws.Range('B42').Value = 'ThetaToken'
$ws.Range('B42').Style = "Normal"

$ws.Range('C42').NumberFormat = "@"
$ws.Range('C42').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('C42').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '2.89'
$ws.Range('D42').Style = "Normal"

$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +7.81%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('B43').NumberFormat = "@"
$ws.Range('B43').Value = 'Stellar'
$ws.Range('B43').Style = "Normal"

$ws.Range('C43').NumberFormat = "@"
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C43').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.142'
$ws.Range('D43').Style = "Normal"

$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +8.54%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +0.41%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +9.98%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.322'
$ws.Range('D46').Style = "Normal"

$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +16.09%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '3.15'
$ws.Range('D47').Style = "Normal"

$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +2.15%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('B48').NumberFormat = "@"
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('B48').Style = "Normal"

$ws.Range('C48').NumberFormat = "@"
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C48').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '2.08'
$ws.Range('D48').Style = "Normal"

$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +6.00%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('B49').NumberFormat = "@"
$ws.Range('B49').Value = 'WEMIXToken'
$ws.Range('B49').Style = "Normal"

$ws.Range('C49').NumberFormat = "@"
$ws.Range('C49').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('C49').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.64'
$ws.Range('D49').Style = "Normal"

$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +6.56%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '142.49'
$ws.Range('D50').Style = "Normal"

$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +0.25%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.82'
$ws.Range('D51').Style = "Normal"

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +6.11%  '
$ws.Range('E51').Style = "Normal"
